$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.805.95"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "2.680.61"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("D5").Value = "'554.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "'156.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("E9").Value = "  -3.29%  "
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("D11").Value = "'5.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("E12").Value = "  -3.44%  "
$ws.Range("D13").Value = "3.156.03"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").Value = "'26.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "62.741.93"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").Value = "2.681.75"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").Value = "'11.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.15%  "
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("D20").Value = "'344.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("E21").Value = "  -5.68%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").Value = "'63.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").Value = "'8.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.86%  "
$ws.Range("D28").Value = "'1.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.50%  "
$ws.Range("D29").Value = "0.0₃0845"
$ws.Range("E29").Value = "  -6.05%  "
$ws.Range("D30").Value = "'7.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "'162.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.96%  "
$ws.Range("D34").Value = "'4.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").Value = "'19.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.99%  "
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("D38").Value = "'339.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("D39").Value = "'6.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").Value = "'0.924"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.88%  "
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("D42").Value = "'38.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").Value = "'20.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.47%  "
$ws.Range("E44").Value = "  -4.36%  "
$ws.Range("D45").Value = "'0.614"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.78%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E47").Value = "  -5.07%  "
$ws.Range("D48").Value = "'10.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'0.0967"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.20%  "
$ws.Range("D50").Value = "'128.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").Value = "'0.0240"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.60%  "
